$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$styleRef = $ws.Range("D5")
$ws.Range("D2").Value = "29.382.86"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.849.46"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D6").Formula = "'0.6287"
$ws.Range("D6").Style = $styleRef.Style
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Formula = "'0.07625"
$ws.Range("D8").Style = $styleRef.Style
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").Formula = "'24.74"
$ws.Range("D10").Style = $styleRef.Style
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("D11").Formula = "'0.07737"
$ws.Range("D11").Style = $styleRef.Style
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D13").Formula = "'0.6792"
$ws.Range("D13").Style = $styleRef.Style
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").Formula = "'0.00001056"
$ws.Range("D14").Style = $styleRef.Style
$ws.Range("E14").Value = "  -3.54%  "
$ws.Range("D15").Formula = "'83.18"
$ws.Range("D15").Style = $styleRef.Style
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").Formula = "'6.167"
$ws.Range("D16").Style = $styleRef.Style
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").Value = "29.398.40"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Formula = "'228.06"
$ws.Range("D18").Style = $styleRef.Style
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Formula = "'12.35"
$ws.Range("D19").Style = $styleRef.Style
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").Formula = "'0.9997"
$ws.Range("D20").Style = $styleRef.Style
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Formula = "'7.494"
$ws.Range("D21").Style = $styleRef.Style
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").Formula = "'1.000"
$ws.Range("D22").Style = $styleRef.Style
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Formula = "'158.80"
$ws.Range("D23").Style = $styleRef.Style
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Formula = "'8.405"
$ws.Range("D25").Style = $styleRef.Style
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").Formula = "'17.70"
$ws.Range("D26").Style = $styleRef.Style
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").Formula = "'1.401"
$ws.Range("D27").Style = $styleRef.Style
$ws.Range("E27").Value = "  +8.01%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Formula = "'0.05609"
$ws.Range("D29").Style = $styleRef.Style
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Formula = "'4.115"
$ws.Range("D30").Style = $styleRef.Style
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D32").Formula = "'1.164"
$ws.Range("D32").Style = $styleRef.Style
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").Formula = "'0.7009"
$ws.Range("D34").Style = $styleRef.Style
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("D35").Formula = "'2.584"
$ws.Range("D35").Style = $styleRef.Style
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "1.234.29"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").Formula = "'0.01806"
$ws.Range("D37").Style = $styleRef.Style
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("D38").Formula = "'2.722"
$ws.Range("D38").Style = $styleRef.Style
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("D39").Formula = "'6.386"
$ws.Range("D39").Style = $styleRef.Style
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("D40").Formula = "'0.9009"
$ws.Range("D40").Style = $styleRef.Style
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Formula = "'101.48"
$ws.Range("D42").Style = $styleRef.Style
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Formula = "'66.02"
$ws.Range("D43").Style = $styleRef.Style
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("D45").Formula = "'0.4001"
$ws.Range("D45").Style = $styleRef.Style
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("E46").Value = "  -5.44%  "
$ws.Range("D47").Formula = "'9.019"
$ws.Range("D47").Style = $styleRef.Style
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").Formula = "'1.679"
$ws.Range("D48").Style = $styleRef.Style
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("D49").Formula = "'0.1132"
$ws.Range("D49").Style = $styleRef.Style
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("D50").Formula = "'0.05706"
$ws.Range("D50").Style = $styleRef.Style
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").Formula = "'0.4627"
$ws.Range("D51").Style = $styleRef.Style
$ws.Range("E51").Value = "  +0.02%  "
